# Fruta / hortaliza, semanal
# Insert a new weekly price observation as row 94 (Vega Monumental Concepción -
# Piña, Segunda, Ecuador), pushing all subsequent rows (old 94-126) down to
# (95-127).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 94..126 down to 95..127, leaving a blank row 94 to fill in.
$ws.Rows.Item(94).Insert()

$newRow = 94
$ws.Cells.Item($newRow, 1).Value = 11
$ws.Cells.Item($newRow, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($newRow, 3).Value = "Bíobío"
$ws.Cells.Item($newRow, 4).Value = 44559
$ws.Cells.Item($newRow, 5).Value = 8
$ws.Cells.Item($newRow, 6).Value = "Fruta"
$ws.Cells.Item($newRow, 7).Value = 100108
$ws.Cells.Item($newRow, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($newRow, 9).Value = 100108005
$ws.Cells.Item($newRow, 10).Value = "Piña"
$ws.Cells.Item($newRow, 11).Value = "Caramelo"
$ws.Cells.Item($newRow, 12).Value = "Segunda"
$ws.Cells.Item($newRow, 13).Value = 200
$ws.Cells.Item($newRow, 14).Value = 15000
$ws.Cells.Item($newRow, 15).Value = 16000
$ws.Cells.Item($newRow, 16).Value = 15500
$ws.Cells.Item($newRow, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item($newRow, 18).Value = "Ecuador"
$ws.Cells.Item($newRow, 19).Value = 1107
$ws.Cells.Item($newRow, 20).Value = 14
